$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")
$ws.Range("C2:C386").Value = 46062
